$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2303.8125
$ws.Range("J100").Value = 5000
$ws.Range("L100").Value = 5000
$ws.Range("N100").Value = -6082
$ws.Range("H135").Value = 1633.1333
$ws.Range("I135").Value = 1678.5
$ws.Range("J135").Value = 998
$ws.Range("K135").Value = 15106.5
$ws.Range("L135").Value = 8982
$ws.Range("M135").Value = -12571.5
$ws.Range("N135").Value = -14052

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3728.3684
$ws.Range("I45").Value = 3112.7778
$ws.Range("K45").Value = 3112.7778
$ws.Range("M45").Value = -2735.7778
$ws.Range("H61").Value = 5962.385
$ws.Range("I61").Value = 8232.714
$ws.Range("K61").Value = 8232.714
$ws.Range("M61").Value = -8020.714
$ws.Range("H74").Value = 1961.625
$ws.Range("I74").Value = 2115.8333
$ws.Range("J74").Value = 1499
$ws.Range("K74").Value = 2115.8333
$ws.Range("L74").Value = 1499
$ws.Range("M74").Value = -1241.8333
$ws.Range("N74").Value = -3247
$ws.Range("H77").Value = 1961.625
$ws.Range("I77").Value = 2115.8333
$ws.Range("J77").Value = 1499
$ws.Range("K77").Value = 10579.1665
$ws.Range("L77").Value = 7495
$ws.Range("M77").Value = -6211.166499999999
$ws.Range("N77").Value = -16231
$ws.Range("H102").Value = 1435.3914
$ws.Range("I102").Value = 1214.9524
$ws.Range("K102").Value = 1214.9524
$ws.Range("M102").Value = 407.0476000000001
$ws.Range("H132").Value = 3144.551
$ws.Range("I132").Value = 2925.1538
$ws.Range("J132").Value = 4000.2
$ws.Range("K132").Value = 8775.4614
$ws.Range("L132").Value = 12000.6
$ws.Range("M132").Value = -6245.4614
$ws.Range("N132").Value = -17060.6
$ws.Range("H136").Value = 5962.385
$ws.Range("I136").Value = 8232.714
$ws.Range("K136").Value = 24698.142
$ws.Range("M136").Value = -22148.142

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2262.6667
$ws.Range("I20").Value = 2144
$ws.Range("J20").Value = 2500
$ws.Range("K20").Value = 2144
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = -1897
$ws.Range("N20").Value = -2994
$ws.Range("H105").Value = 1758.2632
$ws.Range("I105").Value = 1524.0588
$ws.Range("K105").Value = 1524.0588
$ws.Range("M105").Value = 222.9412
$ws.Range("H107").Value = 2076.8064
$ws.Range("I107").Value = 2184.3845
$ws.Range("K107").Value = 2184.3845
$ws.Range("M107").Value = -264.3845000000001
$ws.Range("H134").Value = 3223.7441
$ws.Range("I134").Value = 2753.5588
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8260.6764
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -5725.6764
$ws.Range("N134").Value = -20070

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5261.8335
$ws.Range("I31").Value = 3346.8333
$ws.Range("J31").Value = 5740.5835
$ws.Range("K31").Value = 3346.8333
$ws.Range("L31").Value = 5740.5835
$ws.Range("M31").Value = -3051.8333
$ws.Range("N31").Value = -6330.5835
$ws.Range("H34").Value = 5261.8335
$ws.Range("I34").Value = 3346.8333
$ws.Range("J34").Value = 5740.5835
$ws.Range("K34").Value = 3346.8333
$ws.Range("L34").Value = 5740.5835
$ws.Range("M34").Value = -3144.8333
$ws.Range("N34").Value = -6144.5835
$ws.Range("H58").Value = 3073.4905
$ws.Range("J58").Value = 4399.7334
$ws.Range("L58").Value = 4399.7334
$ws.Range("N58").Value = -4805.7334
$ws.Range("H62").Value = 3914.1
$ws.Range("I62").Value = 2530.3333
$ws.Range("K62").Value = 2530.3333
$ws.Range("M62").Value = -1906.3333
$ws.Range("H65").Value = 3914.1
$ws.Range("I65").Value = 2530.3333
$ws.Range("K65").Value = 12651.6665
$ws.Range("M65").Value = -9531.666499999999
$ws.Range("H132").Value = 3383.8667
$ws.Range("I132").Value = 3196
$ws.Range("K132").Value = 9588
$ws.Range("M132").Value = -7058
$ws.Range("H134").Value = 1994.25
$ws.Range("I134").Value = 1994.091
$ws.Range("K134").Value = 5982.272999999999
$ws.Range("M134").Value = -3447.272999999999
$ws.Range("H136").Value = 3073.4905
$ws.Range("J136").Value = 4399.7334
$ws.Range("L136").Value = 13199.2002
$ws.Range("N136").Value = -18299.2002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9507.571
$ws.Range("I56").Value = 9507.571
$ws.Range("K56").Value = 9507.571
$ws.Range("M56").Value = -8977.571
$ws.Range("H96").Value = 192500
$ws.Range("I96").Value = 200000
$ws.Range("J96").Value = 185000
$ws.Range("K96").Value = 600000
$ws.Range("L96").Value = 555000
$ws.Range("M96").Value = -597941
$ws.Range("N96").Value = -559118
$ws.Range("H113").Value = 1607.75
$ws.Range("I113").Value = 789.7143
$ws.Range("J113").Value = 1944.5883
$ws.Range("K113").Value = 2369.1429
$ws.Range("L113").Value = 5833.7649
$ws.Range("M113").Value = -199.1428999999998
$ws.Range("N113").Value = -10173.7649
$ws.Range("H131").Value = 1544.4865
$ws.Range("I131").Value = 628.4286
$ws.Range("J131").Value = 1758.2333
$ws.Range("K131").Value = 1885.2858
$ws.Range("L131").Value = 5274.699900000001
$ws.Range("M131").Value = 3154.7142
$ws.Range("N131").Value = -15354.6999
$ws.Range("H132").Value = 1668324.1
$ws.Range("I132").Value = 973.5
$ws.Range("J132").Value = 2501999.5
$ws.Range("K132").Value = 8761.5
$ws.Range("L132").Value = 22517995.5
$ws.Range("M132").Value = -6231.5
$ws.Range("N132").Value = -22523055.5
$ws.Range("H137").Value = 15000
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 16857.143
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 50571.429
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -60771.429

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 888.8125
$ws.Range("J97").Value = 2613.6667
$ws.Range("L97").Value = 2613.6667
$ws.Range("N97").Value = -3605.6667
$ws.Range("H107").Value = 777.5625
$ws.Range("I107").Value = 911.8333
$ws.Range("K107").Value = 911.8333
$ws.Range("M107").Value = 1008.1667
$ws.Range("H113").Value = 29829.08
$ws.Range("I113").Value = 14013.6
$ws.Range("J113").Value = 53552.3
$ws.Range("K113").Value = 14013.6
$ws.Range("L113").Value = 53552.3
$ws.Range("M113").Value = -11843.6
$ws.Range("N113").Value = -57892.3

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 20007
$ws.Range("I25").Value = 20007
$ws.Range("K25").Value = 20007
$ws.Range("M25").Value = -19777
$ws.Range("H40").Value = 2153.8948
$ws.Range("I40").Value = 4874.6665
$ws.Range("J40").Value = 1643.75
$ws.Range("K40").Value = 4874.6665
$ws.Range("L40").Value = 1643.75
$ws.Range("M40").Value = -4738.6665
$ws.Range("N40").Value = -1915.75
$ws.Range("H93").Value = 2357.44
$ws.Range("I93").Value = 2173.1428
$ws.Range("J93").Value = 2592
$ws.Range("K93").Value = 2173.1428
$ws.Range("L93").Value = 2592
$ws.Range("M93").Value = -925.1428000000001
$ws.Range("N93").Value = -5088
$ws.Range("H132").Value = 4175.5
$ws.Range("I132").Value = 3333
$ws.Range("K132").Value = 9999
$ws.Range("M132").Value = -7469

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 100000
$ws.Range("J98").Value = 100000
$ws.Range("L98").Value = 100000
$ws.Range("N98").Value = -105990
$ws.Range("H100").Value = 1224.3334
$ws.Range("I100").Value = 436.5
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 873
$ws.Range("L100").Value = 5600
$ws.Range("M100").Value = -332
$ws.Range("N100").Value = -6682
$ws.Range("H107").Value = 761
$ws.Range("I107").Value = 698.3333
$ws.Range("J107").Value = 949
$ws.Range("K107").Value = 2094.9999
$ws.Range("L107").Value = 2847
$ws.Range("M107").Value = -174.9998999999998
$ws.Range("N107").Value = -6687
$ws.Range("H122").Value = 5045.606
$ws.Range("I122").Value = 5431.087
$ws.Range("K122").Value = 16293.261
$ws.Range("M122").Value = -13843.261
